$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix missing style on existing row 33 (I33 should carry the boolean
#     alignment style already used by I2:I32) ---
$ws.Range("I33").Style = $ws.Range("I32").Style

# --- New user records from the 16th May refresh ---
$newRows = @(
    @{ Row = 34; A = 110033; B = 9317596771; C = "Nikola Tesla"; D = "nikola.tesla@xyz.com"; E = 818876434 },
    @{ Row = 35; A = 110034; B = 9317596772; C = "Graham Bell";  D = "graham.bell@xyz.com";  E = 818876435 },
    @{ Row = 36; A = 110035; B = 9317596773; C = "Albert Miles"; D = "albert.miles@xyz.com"; E = 818876436 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = "ACT"
    $ws.Cells.Item($row, 7).Value = "eng"
    $ws.Cells.Item($row, 8).Value = "PWD"
    $ws.Cells.Item($row, 9).Value = $true
    $ws.Cells.Item($row, 10).Value = "superadmin"
    $ws.Cells.Item($row, 11).Value = "now()"
    $ws.Cells.Item($row, 12).Value = "now()"

    # Match the formatting already used in the preceding data rows.
    $ws.Range("D$row").Style = $ws.Range("D32").Style
    $ws.Range("I$row").Style = $ws.Range("I32").Style
}

# --- Reset the stored selection back to the top of the unused columns ---
$ws.Range("M1:XFD1048576").Select()
